{"js": "// Insert Stefan Lyocsa's publications into the \"Main Publication Output\"\n// list, right after item 17 (\"Predicting Retail Customers' ...\") and\n// before the \"Lennart John Baals (ORCID: ...)\" paragraph.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Locate the anchor paragraph - item 17, unique by its text.\nlet anchor = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const t = paragraphs.items[i].text;\n  if (t.indexOf(\"Predicting Retail Customers\") !== -1) {\n    anchor = paragraphs.items[i];\n    break;\n  }\n}\n\nif (!anchor) {\n  throw new Error(\"Anchor paragraph (item 17 / Predicting Retail Customers) not found\");\n}\n\nconst newLines = [\n  \"Stefan Lyocsa (ORCID: 0000-0002-8380-181X):\",\n  \"18. \\\"Macroeconomic environment and the future performance of loans: Evidence from three peer-to-peer platforms\\\". International Review of Financial Analysis. DOI: 10.1016/j.irfa.2024.103416\",\n  \"19. \\\"What drives the uranium sector risk? The role of attention, economic and geopolitical uncertainty\\\". Energy Economics. DOI: 10.1016/j.eneco.2024.107980\",\n  \"20. \\\"Forecasting of clean energy market volatility: The role of oil and the technology sector\\\". Energy Economics. DOI: 10.1016/j.eneco.2024.107451\",\n  \"21. \\\"A Fuzzy Framework for Realized Volatility Prediction\\\" (2025). SSRN.\",\n  \"22. \\\"Alpha-threshold networks in credit risk models\\\" (2025). SSRN.\",\n  \"23. \\\"Do hurricanes cause storm on the stock market?\\\" (2025). SSRN.\"\n];\n\n// Insert each new paragraph right after the previous one, keeping order.\nlet insertAfter = anchor;\nfor (const line of newLines) {\n  insertAfter = insertAfter.insertParagraph(line, \"After\");\n}\n\nawait context.sync();\n", "ps1": "# Insert Stefan Lyocsa's publications into the \"Main Publication Output\"\n# list, right after item 17 (\"Predicting Retail Customers' ...\") and\n# before the \"Lennart John Baals (ORCID: ...)\" paragraph.\n\n$d = $word.ActiveDocument\n\n# Locate the anchor paragraph (item 17) by its unique text.\n$anchorIdx = -1\n$count = $d.Paragraphs.Count\nfor ($i = 1; $i -le $count; $i++) {\n    $t = $d.Paragraphs.Item($i).Range.Text\n    if ($t -like \"*Predicting Retail Customers*\") {\n        $anchorIdx = $i\n        break\n    }\n}\n\nif ($anchorIdx -eq -1) {\n    throw \"Anchor paragraph (item 17 / Predicting Retail Customers) not found\"\n}\n\n$newLines = @(\n    \"Stefan Lyocsa (ORCID: 0000-0002-8380-181X):\",\n    \"18. `\"Macroeconomic environment and the future performance of loans: Evidence from three peer-to-peer platforms`\". International Review of Financial Analysis. DOI: 10.1016/j.irfa.2024.103416\",\n    \"19. `\"What drives the uranium sector risk? The role of attention, economic and geopolitical uncertainty`\". Energy Economics. DOI: 10.1016/j.eneco.2024.107980\",\n    \"20. `\"Forecasting of clean energy market volatility: The role of oil and the technology sector`\". Energy Economics. DOI: 10.1016/j.eneco.2024.107451\",\n    \"21. `\"A Fuzzy Framework for Realized Volatility Prediction`\" (2025). SSRN.\",\n    \"22. `\"Alpha-threshold networks in credit risk models`\" (2025). SSRN.\",\n    \"23. `\"Do hurricanes cause storm on the stock market?`\" (2025). SSRN.\"\n)\n\n$idx = $anchorIdx\nforeach ($line in $newLines) {\n    $cur = $d.Paragraphs.Item($idx).Range\n    $cur.InsertParagraphAfter()\n    $idx = $idx + 1\n    $newPara = $d.Paragraphs.Item($idx).Range\n    $newPara.Text = $line\n}\n"}
